$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 10 (this pushes the existing rows 10-18
# down to rows 12-20, preserving their data/formatting).
$ws.Range("A10:A11").EntireRow.Insert()

# The new rows 10/11 are new weekly entries; seed them from the (now
# shifted) rows 12/13 which hold what used to be rows 10/11, then tweak
# the handful of cells that actually differ for the new week.
$ws.Range("A12:T12").Copy()
$ws.Range("A10:T10").PasteSpecial()

$ws.Range("A13:T13").Copy()
$ws.Range("A11:T11").PasteSpecial()

$excel.CutCopyMode = 0

# Row 10: new date (2021-12-03) and updated price figures for the week.
$ws.Cells.Item(10, 4).Value2 = 44533
$ws.Cells.Item(10, 13).Value2 = 300
$ws.Cells.Item(10, 14).Value2 = 18000
$ws.Cells.Item(10, 15).Value2 = 19000
$ws.Cells.Item(10, 16).Value2 = 18500
$ws.Cells.Item(10, 19).Value2 = 2312

# Row 11: only the date moves to the new week.
$ws.Cells.Item(11, 4).Value2 = 44533
